# Apply cryptocurrency price/volume updates per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.868.37"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "1.630.19"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.67"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2562"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06338"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.48"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07787"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.239"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.634.91"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("D14").Value = "1.854.49"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5521"
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.70"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").Value = "0.0₅7603"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").Value = "25.876.66"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.93"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.415"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.861"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.017"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.893"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.83"
$ws.Range("E26").Value = "  +0.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1251"
$ws.Range("E27").Value = "  +4.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.757"
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.57"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04906"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.235"
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.180"
$ws.Range("E33").Value = "  +0.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.547"
$ws.Range("E34").Value = "  +1.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.372"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.8950"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5522"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.535"
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("D39").Value = "1.115.95"
$ws.Range("E39").Value = "  -2.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01554"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9999"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.570"
$ws.Range("E42").Value = "  +3.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7938"
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.74"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("D45").Value = "1.778.67"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  -12.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4430"
$ws.Range("E47").Value = "  -2.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.71"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05129"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.568"
$ws.Range("E51").Value = "  +3.71%  "
